# Auto-generated edit script applying scheduled-runner price/profit updates
# across the 8 job-leve worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 203.07692
$ws.Range("I4").Value = 194.33333
$ws.Range("J4").Value = 222.75
$ws.Range("K4").Value = 194.33333
$ws.Range("L4").Value = 222.75
$ws.Range("M4").Value = -80.33332999999999
$ws.Range("N4").Value = -450.75
$ws.Range("H18").Value = 465.83334
$ws.Range("J18").Value = 435
$ws.Range("L18").Value = 435
$ws.Range("N18").Value = -1003
$ws.Range("H43").Value = 3248.7144
$ws.Range("J43").Value = 3248.7144
$ws.Range("L43").Value = 3248.7144
$ws.Range("N43").Value = -3386.7144
$ws.Range("H62").Value = 1270.6666
$ws.Range("J62").Value = 1888.4
$ws.Range("L62").Value = 1888.4
$ws.Range("N62").Value = -3136.4
$ws.Range("H65").Value = 1270.6666
$ws.Range("J65").Value = 1888.4
$ws.Range("L65").Value = 9442
$ws.Range("N65").Value = -15682
$ws.Range("H129").Value = 1066.6666
$ws.Range("I129").Value = 436.85715
$ws.Range("J129").Value = 1174.1951
$ws.Range("K129").Value = 1310.57145
$ws.Range("L129").Value = 3522.5853
$ws.Range("M129").Value = 3689.42855
$ws.Range("N129").Value = -13522.5853
$ws.Range("H132").Value = 202479.12
$ws.Range("I132").Value = 224826.38
$ws.Range("J132").Value = 1353.8
$ws.Range("K132").Value = 674479.14
$ws.Range("L132").Value = 4061.4
$ws.Range("M132").Value = -671949.14
$ws.Range("N132").Value = -9121.4
$ws.Range("H137").Value = 3834.7112
$ws.Range("I137").Value = 4430.6665
$ws.Range("J137").Value = 2642.8
$ws.Range("K137").Value = 13291.9995
$ws.Range("L137").Value = 7928.400000000001
$ws.Range("M137").Value = -10741.9995
$ws.Range("N137").Value = -13028.4

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1961172.5
$ws.Range("I2").Value = 391.55554
$ws.Range("J2").Value = 4902344
$ws.Range("K2").Value = 391.55554
$ws.Range("L2").Value = 4902344
$ws.Range("M2").Value = -278.55554
$ws.Range("N2").Value = -4902570
$ws.Range("H3").Value = 37039652
$ws.Range("I3").Value = 125001704
$ws.Range("K3").Value = 125001704
$ws.Range("M3").Value = -125001589
$ws.Range("H61").Value = 2061.0386
$ws.Range("I61").Value = 1828.4474
$ws.Range("J61").Value = 2692.3572
$ws.Range("K61").Value = 1828.4474
$ws.Range("L61").Value = 2692.3572
$ws.Range("M61").Value = -1616.4474
$ws.Range("N61").Value = -3116.3572
$ws.Range("H74").Value = 2225188
$ws.Range("I74").Value = 2857600.2
$ws.Range("J74").Value = 11745.4
$ws.Range("K74").Value = 2857600.2
$ws.Range("L74").Value = 11745.4
$ws.Range("M74").Value = -2856726.2
$ws.Range("N74").Value = -13493.4
$ws.Range("H77").Value = 2225188
$ws.Range("I77").Value = 2857600.2
$ws.Range("J77").Value = 11745.4
$ws.Range("K77").Value = 14288001
$ws.Range("L77").Value = 58727
$ws.Range("M77").Value = -14283633
$ws.Range("N77").Value = -67463
$ws.Range("H112").Value = 30897.666
$ws.Range("J112").Value = 30897.666
$ws.Range("L112").Value = 30897.666
$ws.Range("N112").Value = -33851.666
$ws.Range("H116").Value = 1961172.5
$ws.Range("I116").Value = 391.55554
$ws.Range("J116").Value = 4902344
$ws.Range("K116").Value = 391.55554
$ws.Range("L116").Value = 4902344
$ws.Range("M116").Value = 1902.44446
$ws.Range("N116").Value = -4906932
$ws.Range("H136").Value = 2061.0386
$ws.Range("I136").Value = 1828.4474
$ws.Range("J136").Value = 2692.3572
$ws.Range("K136").Value = 5485.3422
$ws.Range("L136").Value = 8077.071599999999
$ws.Range("M136").Value = -2935.3422
$ws.Range("N136").Value = -13177.0716

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1961172.5
$ws.Range("I3").Value = 391.55554
$ws.Range("J3").Value = 4902344
$ws.Range("K3").Value = 391.55554
$ws.Range("L3").Value = 4902344
$ws.Range("M3").Value = -277.55554
$ws.Range("N3").Value = -4902572

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3710
$ws.Range("I62").Value = 3220
$ws.Range("J62").Value = 4200
$ws.Range("K62").Value = 3220
$ws.Range("L62").Value = 4200
$ws.Range("M62").Value = -2596
$ws.Range("N62").Value = -5448
$ws.Range("H65").Value = 3710
$ws.Range("I65").Value = 3220
$ws.Range("J65").Value = 4200
$ws.Range("K65").Value = 16100
$ws.Range("L65").Value = 21000
$ws.Range("M65").Value = -12980
$ws.Range("N65").Value = -27240

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1226.3636
$ws.Range("I3").Value = 832.2222
$ws.Range("K3").Value = 2496.6666
$ws.Range("M3").Value = -2384.6666
$ws.Range("H4").Value = 799.875
$ws.Range("J4").Value = 1351
$ws.Range("L4").Value = 4053
$ws.Range("N4").Value = -4277
$ws.Range("H63").Value = 3318.3333
$ws.Range("I63").Value = 955
$ws.Range("K63").Value = 2865
$ws.Range("M63").Value = -2116
$ws.Range("H66").Value = 3318.3333
$ws.Range("I66").Value = 955
$ws.Range("K66").Value = 8595
$ws.Range("M66").Value = -4851
$ws.Range("H68").Value = 858.9136
$ws.Range("I68").Value = 617.13043
$ws.Range("J68").Value = 1176.6857
$ws.Range("K68").Value = 1851.39129
$ws.Range("L68").Value = 3530.0571
$ws.Range("M68").Value = -1040.39129
$ws.Range("N68").Value = -5152.0571
$ws.Range("H71").Value = 858.9136
$ws.Range("I71").Value = 617.13043
$ws.Range("J71").Value = 1176.6857
$ws.Range("K71").Value = 5554.173870000001
$ws.Range("L71").Value = 10590.1713
$ws.Range("M71").Value = -1498.173870000001
$ws.Range("N71").Value = -18702.1713
$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -19118
$ws.Range("H129").Value = 859.4666999999999
$ws.Range("I129").Value = 312.7143
$ws.Range("J129").Value = 1337.875
$ws.Range("K129").Value = 938.1428999999999
$ws.Range("L129").Value = 4013.625
$ws.Range("M129").Value = 4061.8571
$ws.Range("N129").Value = -14013.625
$ws.Range("H131").Value = 1151073.1
$ws.Range("J131").Value = 1430405.5
$ws.Range("L131").Value = 4291216.5
$ws.Range("N131").Value = -4301296.5
$ws.Range("H133").Value = 7018.737
$ws.Range("H134").Value = 29414276
$ws.Range("I134").Value = 29414276
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 88242828
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -88237758
$ws.Range("N134").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67.5
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 75
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 38
$ws.Range("N2").Value = -271
$ws.Range("H3").Value = 1556.6522
$ws.Range("I3").Value = 718.3182
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 718.3182
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = -602.3182
$ws.Range("N3").Value = -20232
$ws.Range("H10").Value = 112922.22
$ws.Range("I10").Value = 144542.86
$ws.Range("J10").Value = 2250
$ws.Range("K10").Value = 144542.86
$ws.Range("L10").Value = 2250
$ws.Range("M10").Value = -144373.86
$ws.Range("N10").Value = -2588
$ws.Range("H126").Value = 1262.5
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -830
$ws.Range("N126").Value = -10190

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 33338334
$ws.Range("J2").Value = 50004000
$ws.Range("L2").Value = 50004000
$ws.Range("N2").Value = -50004224
$ws.Range("H110").Value = 23000
$ws.Range("J110").Value = 23000
$ws.Range("L110").Value = 23000
$ws.Range("N110").Value = -31180

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 580
$ws.Range("J5").Value = 580
$ws.Range("L5").Value = 580
$ws.Range("N5").Value = -804
$ws.Range("H6").Value = 2000
$ws.Range("J6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("N6").Value = -2230
$ws.Range("H9").Value = 19500
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 10725
$ws.Range("J12").Value = 10725
$ws.Range("L12").Value = 10725
$ws.Range("N12").Value = -11009
$ws.Range("H126").Value = 918.53845
$ws.Range("I126").Value = 893.1
$ws.Range("J126").Value = 1003.3333
$ws.Range("K126").Value = 2679.3
$ws.Range("L126").Value = 3009.9999
$ws.Range("M126").Value = -209.3000000000002
$ws.Range("N126").Value = -7949.9999
